$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 49208
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 44005.25
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 44005.25
$ws.Range("M21").Value = -69551
$ws.Range("N21").Value = -44941.25
$ws.Range("H23").Value = 49208
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 44005.25
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 44005.25
$ws.Range("M23").Value = -69785
$ws.Range("N23").Value = -44473.25
$ws.Range("H34").Value = 15647.728
$ws.Range("I34").Value = 1039.7142
$ws.Range("K34").Value = 1039.7142
$ws.Range("M34").Value = -836.7141999999999
$ws.Range("H36").Value = 15647.728
$ws.Range("I36").Value = 1039.7142
$ws.Range("K36").Value = 1039.7142
$ws.Range("M36").Value = -324.7141999999999
$ws.Range("H43").Value = 15996963
$ws.Range("I43").Value = 35739540
$ws.Range("J43").Value = 202900.4
$ws.Range("K43").Value = 35739540
$ws.Range("L43").Value = 202900.4
$ws.Range("M43").Value = -35739471
$ws.Range("N43").Value = -203038.4
$ws.Range("H137").Value = 1020.2923
$ws.Range("I137").Value = 820.38
$ws.Range("J137").Value = 1686.6666
$ws.Range("K137").Value = 2461.14
$ws.Range("L137").Value = 5059.9998
$ws.Range("M137").Value = 88.86000000000013
$ws.Range("N137").Value = -10159.9998
$ws.Range("H138").Value = 1801.67
$ws.Range("I138").Value = 805.8214
$ws.Range("J138").Value = 2188.9443
$ws.Range("K138").Value = 2417.4642
$ws.Range("L138").Value = 6566.8329
$ws.Range("M138").Value = 2722.5358
$ws.Range("N138").Value = -16846.8329

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 14872.857
$ws.Range("I26").Value = 1268
$ws.Range("J26").Value = 25076.5
$ws.Range("K26").Value = 1268
$ws.Range("L26").Value = 25076.5
$ws.Range("M26").Value = -938
$ws.Range("N26").Value = -25736.5
$ws.Range("H32").Value = 1235040.6
$ws.Range("I32").Value = 1463517.4
$ws.Range("J32").Value = 77425.13
$ws.Range("K32").Value = 1463517.4
$ws.Range("L32").Value = 77425.13
$ws.Range("M32").Value = -1463230.4
$ws.Range("N32").Value = -77999.13
$ws.Range("H39").Value = 22288
$ws.Range("I39").Value = 3196
$ws.Range("J39").Value = 70018
$ws.Range("K39").Value = 3196
$ws.Range("L39").Value = 70018
$ws.Range("M39").Value = -2676
$ws.Range("N39").Value = -71058
$ws.Range("H61").Value = 1625.5476
$ws.Range("I61").Value = 1485.1143
$ws.Range("J61").Value = 2327.7144
$ws.Range("K61").Value = 1485.1143
$ws.Range("L61").Value = 2327.7144
$ws.Range("M61").Value = -1273.1143
$ws.Range("N61").Value = -2751.7144
$ws.Range("H136").Value = 1625.5476
$ws.Range("I136").Value = 1485.1143
$ws.Range("J136").Value = 2327.7144
$ws.Range("K136").Value = 4455.3429
$ws.Range("L136").Value = 6983.1432
$ws.Range("M136").Value = -1905.3429
$ws.Range("N136").Value = -12083.1432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 800.84
$ws.Range("I20").Value = 687.8823
$ws.Range("J20").Value = 1040.875
$ws.Range("K20").Value = 687.8823
$ws.Range("L20").Value = 1040.875
$ws.Range("M20").Value = -440.8823
$ws.Range("N20").Value = -1534.875
$ws.Range("H62").Value = 70000
$ws.Range("J62").Value = 70000
$ws.Range("L62").Value = 70000
$ws.Range("N62").Value = -71372
$ws.Range("H65").Value = 70000
$ws.Range("J65").Value = 70000
$ws.Range("L65").Value = 210000
$ws.Range("N65").Value = -216864

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1889287.8
$ws.Range("I6").Value = 5003100
$ws.Range("J6").Value = 21000.4
$ws.Range("K6").Value = 5003100
$ws.Range("L6").Value = 21000.4
$ws.Range("M6").Value = -5002987
$ws.Range("N6").Value = -21226.4
$ws.Range("H19").Value = 26817.875
$ws.Range("I19").Value = 905
$ws.Range("K19").Value = 905
$ws.Range("M19").Value = -735
$ws.Range("H24").Value = 26817.875
$ws.Range("I24").Value = 905
$ws.Range("K24").Value = 905
$ws.Range("M24").Value = -735
$ws.Range("H25").Value = 31357
$ws.Range("J25").Value = 35622.285
$ws.Range("L25").Value = 35622.285
$ws.Range("N25").Value = -35970.285
$ws.Range("H58").Value = 3864.6
$ws.Range("I58").Value = 1383
$ws.Range("J58").Value = 5519
$ws.Range("K58").Value = 1383
$ws.Range("L58").Value = 5519
$ws.Range("M58").Value = -1180
$ws.Range("N58").Value = -5925
$ws.Range("H136").Value = 3864.6
$ws.Range("I136").Value = 1383
$ws.Range("J136").Value = 5519
$ws.Range("K136").Value = 4149
$ws.Range("L136").Value = 16557
$ws.Range("M136").Value = -1599
$ws.Range("N136").Value = -21657

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 756.25
$ws.Range("J17").Value = 800
$ws.Range("L17").Value = 2400
$ws.Range("N17").Value = -2738
$ws.Range("H32").Value = 40638.617
$ws.Range("I32").Value = 255001
$ws.Range("J32").Value = 1663.6364
$ws.Range("K32").Value = 765003
$ws.Range("L32").Value = 4990.9092
$ws.Range("M32").Value = -764720
$ws.Range("N32").Value = -5556.9092
$ws.Range("H46").Value = 2140.8518
$ws.Range("I46").Value = 633.3333
$ws.Range("J46").Value = 2571.5715
$ws.Range("K46").Value = 1899.9999
$ws.Range("L46").Value = 7714.7145
$ws.Range("M46").Value = -1808.9999
$ws.Range("N46").Value = -7896.7145
$ws.Range("H113").Value = 763.7308
$ws.Range("I113").Value = 374.41177
$ws.Range("J113").Value = 1499.1111
$ws.Range("K113").Value = 1123.23531
$ws.Range("L113").Value = 4497.3333
$ws.Range("M113").Value = 1046.76469
$ws.Range("N113").Value = -8837.3333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 23702
$ws.Range("I22").Value = 5354
$ws.Range("J22").Value = 29818
$ws.Range("K22").Value = 5354
$ws.Range("L22").Value = 29818
$ws.Range("M22").Value = -4825
$ws.Range("N22").Value = -30876
$ws.Range("H70").Value = 4096.143
$ws.Range("I70").Value = 4013.3333
$ws.Range("J70").Value = 4158.25
$ws.Range("K70").Value = 4013.3333
$ws.Range("L70").Value = 4158.25
$ws.Range("M70").Value = -3743.3333
$ws.Range("N70").Value = -4698.25
$ws.Range("H73").Value = 4096.143
$ws.Range("I73").Value = 4013.3333
$ws.Range("J73").Value = 4158.25
$ws.Range("K73").Value = 4013.3333
$ws.Range("L73").Value = 4158.25
$ws.Range("M73").Value = -3077.3333
$ws.Range("N73").Value = -6030.25
$ws.Range("H97").Value = 963.3570999999999
$ws.Range("I97").Value = 922.0833
$ws.Range("J97").Value = 1211
$ws.Range("K97").Value = 922.0833
$ws.Range("L97").Value = 1211
$ws.Range("M97").Value = -426.0833
$ws.Range("N97").Value = -2203

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1700
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1405
$ws.Range("H27").Value = 1700
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1593
$ws.Range("H132").Value = 212253.12
$ws.Range("I132").Value = 52390.848
$ws.Range("J132").Value = 558621.4
$ws.Range("K132").Value = 157172.544
$ws.Range("L132").Value = 1675864.2
$ws.Range("M132").Value = -154642.544
$ws.Range("N132").Value = -1680924.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1818.8823
$ws.Range("I96").Value = 1373.1
$ws.Range("J96").Value = 2455.7144
$ws.Range("K96").Value = 1373.1
$ws.Range("L96").Value = 2455.7144
$ws.Range("M96").Value = -0.09999999999990905
$ws.Range("N96").Value = -5201.7144
$ws.Range("H107").Value = 1283
$ws.Range("I107").Value = 1116.2
$ws.Range("J107").Value = 1700
$ws.Range("K107").Value = 3348.6
$ws.Range("L107").Value = 5100
$ws.Range("M107").Value = -1428.6
$ws.Range("N107").Value = -8940
